# Generate Report for Handoff
#
# For the "82429938-e6c9-4c64-aeed-848f6d261f77" source file row (row 5) on
# both the "zh-cn" and "de-de" localization status sheets, stamp the
# "Latest Handoff Datetime" column (D) with the timestamp of the handoff
# that was just generated.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-25 03:37:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-25 03:38:07"
